$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.23780120935333982
$ws.Range("C2").Value = 0.0070119080416769395
$ws.Range("D2").Value = -0.2515443282740501
$ws.Range("E2").Value = -0.22405809043262956
$ws.Range("B3").Value = 0.11216581351611048
$ws.Range("C3").Value = 0.0026282775840334692
$ws.Range("D3").Value = 0.10701447222935702
$ws.Range("E3").Value = 0.11731715480286393
$ws.Range("B4").Value = 0.20032632630161099
$ws.Range("C4").Value = 0.0092818499759300538
$ws.Range("D4").Value = 0.1821341926818878
$ws.Range("E4").Value = 0.21851845992133417
$ws.Range("B5").Value = -0.089584143710899664
$ws.Range("C5").Value = 0.0037285983220680993
$ws.Range("D5").Value = -0.096892079353877753
$ws.Range("E5").Value = -0.082276208067921575
$ws.Range("B6").Value = 0.068711442454373065
$ws.Range("C6").Value = 0.0022530516772297834
$ws.Range("D6").Value = 0.064295531906956804
$ws.Range("E6").Value = 0.073127353001789325
$ws.Range("B7").Value = 0.068241281570436685
$ws.Range("C7").Value = 0.0052427943682260003
$ws.Range("D7").Value = 0.057965569218800909
$ws.Range("E7").Value = 0.078516993922072453
$ws.Range("B8").Value = -0.26944287625901975
$ws.Range("C8").Value = 0.0058246516482187097
$ws.Range("D8").Value = -0.28085901004134994
$ws.Range("E8").Value = -0.25802674247668955
$ws.Range("B9").Value = 0.21849136632133936
$ws.Range("C9").Value = 0.0026502569184737508
$ws.Range("D9").Value = 0.213296946231325
$ws.Range("E9").Value = 0.22368578641135373
$ws.Range("B10").Value = 0.16993448992396976
$ws.Range("C10").Value = 0.0087939733302174906
$ws.Range("D10").Value = 0.15269857916410334
$ws.Range("E10").Value = 0.18717040068383617
$ws.Range("B11").Value = -0.11752637150882232
$ws.Range("C11").Value = 0.0036356244275116715
$ws.Range("D11").Value = -0.12465208123760793
$ws.Range("E11").Value = -0.11040066178003671
$ws.Range("B12").Value = 0.12101164159224552
$ws.Range("C12").Value = 0.0025790760690650185
$ws.Range("D12").Value = 0.1159567334731489
$ws.Range("E12").Value = 0.12606654971134215
$ws.Range("B13").Value = 0.077874113396346298
$ws.Range("C13").Value = 0.0057576624560337387
$ws.Range("D13").Value = 0.066589275758124233
$ws.Range("E13").Value = 0.089158951034568362
$ws.Range("B14").Value = -0.24419380014280956
$ws.Range("C14").Value = 0.0048031904849626713
$ws.Range("D14").Value = -0.2536079022162131
$ws.Range("E14").Value = -0.23477969806940599
$ws.Range("B15").Value = 0.28270976582374718
$ws.Range("C15").Value = 0.0026393387243403362
$ws.Range("D15").Value = 0.27753674505036435
$ws.Range("E15").Value = 0.28788278659713001
$ws.Range("B16").Value = 0.085910302406704511
$ws.Range("C16").Value = 0.0081868862611263384
$ws.Range("D16").Value = 0.069864263181968755
$ws.Range("E16").Value = 0.10195634163144027
$ws.Range("B17").Value = -0.12178834401309085
$ws.Range("C17").Value = 0.003396818525197192
$ws.Range("D17").Value = -0.12844600167122339
$ws.Range("E17").Value = -0.11513068635495831
$ws.Range("B18").Value = 0.16659143793455386
$ws.Range("C18").Value = 0.0027103059322577706
$ws.Range("D18").Value = 0.16127932340387532
$ws.Range("E18").Value = 0.1719035524652324
$ws.Range("B19").Value = 0.053743469356032587
$ws.Range("C19").Value = 0.0059615228504391561
$ws.Range("D19").Value = 0.042059071745461286
$ws.Range("E19").Value = 0.065427866966603881
$ws.Range("B20").Value = -0.20487305481564672
$ws.Range("C20").Value = 0.0040805792078984484
$ws.Range("D20").Value = -0.21287086154474269
$ws.Range("E20").Value = -0.19687524808655074
$ws.Range("B21").Value = 0.30555608605670598
$ws.Range("C21").Value = 0.0026877674849701736
$ws.Range("D21").Value = 0.30028814643775892
$ws.Range("E21").Value = 0.31082402567565304
$ws.Range("B22").Value = 0.021694275673353813
$ws.Range("C22").Value = 0.0077184211339620537
$ws.Range("D22").Value = 0.0065664133434885572
$ws.Range("E22").Value = 0.036822138003219068
$ws.Range("B23").Value = -0.1143508439081759
$ws.Range("C23").Value = 0.0032660070989607667
$ws.Range("D23").Value = -0.12075211527802311
$ws.Range("E23").Value = -0.10794957253832869
$ws.Range("B24").Value = 0.19926924539985302
$ws.Range("C24").Value = 0.0029044336007393286
$ws.Range("D24").Value = 0.19357664673405378
$ws.Range("E24").Value = 0.20496184406565227
$ws.Range("B25").Value = 0.042788548428895068
$ws.Range("C25").Value = 0.0061878274014262862
$ws.Range("D25").Value = 0.030660601003762437
$ws.Range("E25").Value = 0.054916495854027698
$ws.Range("B26").Value = -0.17115520221954347
$ws.Range("C26").Value = 0.003720958722136835
$ws.Range("D26").Value = -0.17844816412284173
$ws.Range("E26").Value = -0.16386224031624522
$ws.Range("B27").Value = 0.30652632479985342
$ws.Range("C27").Value = 0.0028732812569208913
$ws.Range("D27").Value = 0.30089478403066389
$ws.Range("E27").Value = 0.31215786556904296
$ws.Range("B28").Value = -0.028759102420072236
$ws.Range("C28").Value = 0.0078007699643280774
$ws.Range("D28").Value = -0.044048365863867532
$ws.Range("E28").Value = -0.013469838976276941
$ws.Range("B29").Value = -0.10942786563243886
$ws.Range("C29").Value = 0.0031164897527257117
$ws.Range("D29").Value = -0.11553608769812058
$ws.Range("E29").Value = -0.10331964356675714
$ws.Range("B30").Value = 0.21057787018812663
$ws.Range("C30").Value = 0.0030557086130297492
$ws.Range("D30").Value = 0.20458877724787952
$ws.Range("E30").Value = 0.21656696312837373
$ws.Range("B31").Value = 0.032960362500908003
$ws.Range("C31").Value = 0.0064263377769409186
$ws.Range("D31").Value = 0.020364942228371232
$ws.Range("E31").Value = 0.04555578277344477
$ws.Range("B32").Value = -0.13817285652500857
$ws.Range("C32").Value = 0.0034683280441769131
$ws.Range("D32").Value = -0.14497067025614219
$ws.Range("E32").Value = -0.13137504279387496
$ws.Range("B33").Value = 0.27323389489622568
$ws.Range("C33").Value = 0.0030266646791453598
$ws.Range("D33").Value = 0.26730172745030745
$ws.Range("E33").Value = 0.27916606234214392
$ws.Range("B34").Value = -0.064929843386910976
$ws.Range("C34").Value = 0.007798370925615518
$ws.Range("D34").Value = -0.080214404790387703
$ws.Range("E34").Value = -0.049645281983434257
$ws.Range("B35").Value = -0.09825570228743026
$ws.Range("C35").Value = 0.0030590246679594976
$ws.Range("D35").Value = -0.10425129459122419
$ws.Range("E35").Value = -0.092260109983636332
$ws.Range("B36").Value = 0.20955131050005238
$ws.Range("C36").Value = 0.0032846075243989437
$ws.Range("D36").Value = 0.20311358288035125
$ws.Range("E36").Value = 0.21598903811975351
$ws.Range("B37").Value = 0.030426522718642064
$ws.Range("C37").Value = 0.0068686521207487394
$ws.Range("D37").Value = 0.016964180219760257
$ws.Range("E37").Value = 0.043888865217523873
$ws.Range("B38").Value = -0.10818951530566451
$ws.Range("C38").Value = 0.0030917238565077013
$ws.Range("D38").Value = -0.11424919669016546
$ws.Range("E38").Value = -0.10212983392116355
$ws.Range("B39").Value = 0.18707083643902572
$ws.Range("C39").Value = 0.0029554530354155913
$ws.Range("D39").Value = 0.18127824157199768
$ws.Range("E39").Value = 0.19286343130605377
$ws.Range("B40").Value = -0.07507097282684741
$ws.Range("C40").Value = 0.007118902377888362
$ws.Range("D40").Value = -0.089023797276731209
$ws.Range("E40").Value = -0.061118148376963617
$ws.Range("B41").Value = -0.089586449774134255
$ws.Range("C41").Value = 0.0031558678574501054
$ws.Range("D41").Value = -0.095771851688705861
$ws.Range("E41").Value = -0.083401047859562649
$ws.Range("B42").Value = 0.19907249746531222
$ws.Range("C42").Value = 0.0036711333802549326
$ws.Range("D42").Value = 0.19187719130403844
$ws.Range("E42").Value = 0.206267803626586
$ws.Range("B43").Value = 0.02096221083603202
$ws.Range("C43").Value = 0.0077484634780843279
$ws.Range("D43").Value = 0.0057754657005574233
$ws.Range("E43").Value = 0.036148955971506615
$ws.Range("B44").Value = -0.099649642653548892
$ws.Range("C44").Value = 0.0032959646931615149
$ws.Range("D44").Value = -0.10610962964529899
$ws.Range("E44").Value = -0.093189655661798793
$ws.Range("B45").Value = 0.081373064213153878
$ws.Range("C45").Value = 0.0032206477266531314
$ws.Range("D45").Value = 0.075060696103642327
$ws.Range("E45").Value = 0.08768543232266543
$ws.Range("B46").Value = -0.059610815238342135
$ws.Range("C46").Value = 0.0071379534919600455
$ws.Range("D46").Value = -0.07360097927178906
$ws.Range("E46").Value = -0.04562065120489521
$ws.Range("B47").Value = -0.077722117631489382
$ws.Range("C47").Value = 0.0033790285346019851
$ws.Range("D47").Value = -0.084344907466614688
$ws.Range("E47").Value = -0.071099327796364076
$ws.Range("B48").Value = 0.16393014435956857
$ws.Range("C48").Value = 0.0042164964518062176
$ws.Range("D48").Value = 0.1556659437010342
$ws.Range("E48").Value = 0.17219434501810293
$ws.Range("B49").Value = 0.0060801871156313777
$ws.Range("C49").Value = 0.0089248794825567603
$ws.Range("D49").Value = -0.011412296452210872
$ws.Range("E49").Value = 0.023572670683473625
$ws.Range("B50").Value = -0.094680740857441298
$ws.Range("C50").Value = 0.0039287815135156849
$ws.Range("D50").Value = -0.10238102888643563
$ws.Range("E50").Value = -0.086980452828446964
$ws.Range("B51").Value = -0.011733922865044511
$ws.Range("C51").Value = 0.0036514784695653563
$ws.Range("D51").Value = -0.018890705661592466
$ws.Range("E51").Value = -0.0045771400684965552
$ws.Range("B52").Value = -0.00014670012356988375
$ws.Range("C52").Value = 0.0077763052423769569
$ws.Range("D52").Value = -0.015388013482860653
$ws.Range("E52").Value = 0.015094613235720887
$ws.Range("B53").Value = -0.063111966813637427
$ws.Range("C53").Value = 0.0040657144072792795
$ws.Range("D53").Value = -0.071080639399062451
$ws.Range("E53").Value = -0.05514329422821241
$ws.Range("B54").Value = 0.10621357032378968
$ws.Range("C54").Value = 0.0053890156677009653
$ws.Range("D54").Value = 0.095651268816150922
$ws.Range("E54").Value = 0.11677587183142844
$ws.Range("B55").Value = -0.01036103988011374
$ws.Range("C55").Value = 0.010980266934485139
$ws.Range("D55").Value = -0.031882018319923315
$ws.Range("E55").Value = 0.011159938559695838
